$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101:150 down to 102:151
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record
$ws.Range("A101").Value = 9
$ws.Range("B101").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C101").Value = "Metropolitana"
$ws.Range("D101").Value = 44460
$ws.Range("E101").Value = 13
$ws.Range("F101").Value = 300000001
$ws.Range("G101").Value = "Rabanito"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 5200
$ws.Range("K101").Value = 3500
$ws.Range("L101").Value = 4000
$ws.Range("M101").Value = 3750
$ws.Range("N101").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O101").Value = "Provincia de Chacabuco"
$ws.Range("P101").Value = 38
$ws.Range("Q101").Value = 100
$ws.Range("R101").Value = "Hortaliza"
